# Fix Title, Subtitle, Headings in source document.
#
#  - Title      -> Heading 1   (keeps centered alignment)
#  - Subtitle   -> (no style / Normal), keeps centered alignment,
#                  and its Google-Docs bookmark is dropped
#  - The four top-level "Heading 1" section headers (Summary,
#    Prerequisites, Effort, Description) become "Heading 2"
#
$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Title paragraph -> Heading 1, keep centered alignment ----------
$pTitle = $d.Paragraphs.Item(1)
$pTitle.Range.InsertXML(@"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
    <w:jc w:val="center"/>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_gjdgxs" w:colFirst="0" w:colLast="0"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:rtl w:val="0"/>
    </w:rPr>
    <w:t xml:space="preserve">ACCESS Affinity Groups</w:t>
  </w:r>
</w:p>
"@) | Out-Null

# --- 2. Subtitle paragraph -> drop the style + its bookmark -------------
$pSubtitle = $d.Paragraphs.Item(2)
$pSubtitle.Range.InsertXML(@"
<w:p $wNs>
  <w:pPr>
    <w:jc w:val="center"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rtl w:val="0"/>
    </w:rPr>
    <w:t xml:space="preserve">Version 1 - 2022/09/26</w:t>
  </w:r>
</w:p>
"@) | Out-Null

# --- 3. Promote the four Heading 1 section titles to Heading 2 ---------
# (bookmark names are unchanged; Word renumbers the bookmark ids itself)
$headingTexts = @("Summary", "Prerequisites", "Effort", "Description")
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if (($p.Style.NameLocal -eq "Heading 1") -and ($headingTexts -contains $t)) {
        $p.Style = "Heading 2"
    }
}
